# Daily attendance processing - 2026-01-06 19:56:21
#
# The "Recorded By" column (G) lists the users who touched each attendance
# record, as a comma-separated string. Re-normalise the ordering for every
# row whose list still has the stale ordering coming out of the overnight
# ingest job: move the LAST name in the list to the FRONT (a right-rotation
# of the comma-separated tokens), leaving already-normalised / single-name
# rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$rowCount = $usedRange.Rows.Count
$lastRow = $firstRow + $rowCount - 1

# Values that are known to still carry the pre-normalisation ordering.
$staleValues = @(
    "System, backup@backdoor.com, system",
    "System, dnasr281@gmail.com"
)

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($null -eq $current) {
        continue
    }

    $currentText = [string]$current

    if ($staleValues -contains $currentText) {
        $tokens = $currentText.Split(",")
        $trimmed = @()
        foreach ($t in $tokens) {
            $trimmed += $t.Trim()
        }

        if ($trimmed.Count -gt 1) {
            $lastToken = $trimmed[$trimmed.Count - 1]
            $rest = $trimmed[0..($trimmed.Count - 2)]
            $rotated = @($lastToken) + $rest
            $newValue = [string]::Join(", ", $rotated)
            $cell.Value = $newValue
        }
    }
}
